$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (date series) and append the next 6 days.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$dates = @("09-10-2021", "10-10-2021", "11-10-2021", "12-10-2021", "13-10-2021", "14-10-2021")

foreach ($d in $dates) {
    $lastRow = $lastRow + 1

    $cellA = $ws.Cells.Item($lastRow, 1)
    # Enter the date string as a text formula, then copy/paste-as-values
    # over itself. This converts it to a literal (shared-string) text
    # value "dd-mm-yyyy" without Excel's autodetection turning it into a
    # real date serial, and without mutating the cell's number format/style.
    $cellA.Formula = '="' + $d + '"'
    $cellA.Copy()
    $cellA.PasteSpecial(-4163)

    $ws.Cells.Item($lastRow, 2).Value = 3623
    $ws.Cells.Item($lastRow, 3).Value = 240
}

$excel.CutCopyMode = 0
